# Trade #10 closed at 2026-02-17 23:53:04 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.29
$summary.Range("B4").Value = 0.29
$summary.Range("B6").Value = 10
$summary.Range("B7").Value = 6
$summary.Range("B9").Value = 60

# ---- Strategy Status sheet (MarketMaking row) ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.29
$status.Range("D6").Value = 10
$status.Range("E6").Value = 0.29
$status.Range("F6").Value = 0.29
$status.Range("G6").Value = 60

# ---- New trade row data (Trade #10) ----
$newRow = @(10, "2026-02-17", "23:52:58", "MarketMaking", "UP", 0.9, 0.93, "CLOSED", 3.3333, 0.03, 100.29, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

function Write-TradeRow($sheet, $rowNum, $values) {
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $sheet.Cells.Item($rowNum, $col)
        if ($col -eq 2) {
            # Date column ("2026-02-17") - force text so it is not
            # auto-coerced into a date serial number.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $values[$col - 1]
    }
}

# ---- All Trades sheet: append new row 11 ----
$allTrades = $wb.Worksheets.Item("All Trades")
Write-TradeRow $allTrades 11 $newRow

# ---- MarketMaking sheet: append new row 11 ----
$mm = $wb.Worksheets.Item("MarketMaking")
Write-TradeRow $mm 11 $newRow
